# "Adjustments to publication notice and code-check"
#
# The code-check worksheet's column-D header reads "Replicated?" in three
# places (the row-1 heading and the two repeated section headers in rows 6
# and 12); all three cells point at the same shared-string entry. Replacing
# it once, workbook-wide, updates every occurrence and keeps them sharing a
# single (now-renamed) string, exactly as the authored edit does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("Replicated?", "Reproduced?") | Out-Null

# The page setup no longer pins a specific first page number - restore the
# "automatic" page-numbering default (this drops the explicit
# firstPageNumber/useFirstPageNumber overrides on save).
$ws.PageSetup.FirstPageNumber = -4105

# Finally, leave the active selection on D12, matching the cursor position
# recorded in the saved workbook.
$ws.Range("D12").Select() | Out-Null
